$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1880.2222
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 2074.5715
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 2074.5715
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -2424.5715

$ws.Range("H58").Value = 1582.1428
$ws.Range("J58").Value = 1761.5385
$ws.Range("L58").Value = 5284.6155
$ws.Range("N58").Value = -5584.6155

$ws.Range("H64").Value = 3046.6667
$ws.Range("I64").Value = 2922.5806
$ws.Range("J64").Value = 3321.4285
$ws.Range("K64").Value = 2922.5806
$ws.Range("L64").Value = 3321.4285
$ws.Range("M64").Value = -2674.5806
$ws.Range("N64").Value = -3817.4285

$ws.Range("H67").Value = 3046.6667
$ws.Range("I67").Value = 2922.5806
$ws.Range("J67").Value = 3321.4285
$ws.Range("K67").Value = 2922.5806
$ws.Range("L67").Value = 3321.4285
$ws.Range("M67").Value = -2064.5806
$ws.Range("N67").Value = -5037.4285

$ws.Range("H76").Value = 3019.6
$ws.Range("I76").Value = 3010.6667
$ws.Range("K76").Value = 3010.6667
$ws.Range("M76").Value = -2695.6667

$ws.Range("H79").Value = 3019.6
$ws.Range("I79").Value = 3010.6667
$ws.Range("K79").Value = 3010.6667
$ws.Range("M79").Value = -1918.6667

$ws.Range("H99").Value = 812.06665
$ws.Range("I99").Value = 812.06665
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2436.19995
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -938.1999500000002
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3081.77
$ws.Range("I32").Value = 3072.4949
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 3072.4949
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -2785.4949
$ws.Range("N32").Value = -4574

$ws.Range("H88").Value = 2046.6154
$ws.Range("I88").Value = 1966.6666
$ws.Range("J88").Value = 2070.6
$ws.Range("K88").Value = 1966.6666
$ws.Range("L88").Value = 2070.6
$ws.Range("M88").Value = -1560.6666
$ws.Range("N88").Value = -2882.6

$ws.Range("H91").Value = 2046.6154
$ws.Range("I91").Value = 1966.6666
$ws.Range("J91").Value = 2070.6
$ws.Range("K91").Value = 1966.6666
$ws.Range("L91").Value = 2070.6
$ws.Range("M91").Value = -562.6666
$ws.Range("N91").Value = -4878.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2445
$ws.Range("I105").Value = 2408.3333
$ws.Range("K105").Value = 2408.3333
$ws.Range("M105").Value = -661.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2422.7778
$ws.Range("I62").Value = 2434.1667
$ws.Range("K62").Value = 2434.1667
$ws.Range("M62").Value = -1810.1667

$ws.Range("H65").Value = 2422.7778
$ws.Range("I65").Value = 2434.1667
$ws.Range("K65").Value = 12170.8335
$ws.Range("M65").Value = -9050.833500000001

$ws.Range("H74").Value = 23333.334
$ws.Range("J74").Value = 23333.334
$ws.Range("L74").Value = 23333.334
$ws.Range("N74").Value = -25081.334

$ws.Range("H77").Value = 23333.334
$ws.Range("J77").Value = 23333.334
$ws.Range("L77").Value = 70000.00199999999
$ws.Range("N77").Value = -78736.00199999999

$ws.Range("H141").Value = 51089.1
$ws.Range("J141").Value = 51089.1
$ws.Range("L141").Value = 51089.1
$ws.Range("N141").Value = -61449.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4507.273
$ws.Range("I3").Value = 3082.8572
$ws.Range("J3").Value = 7000
$ws.Range("K3").Value = 9248.571599999999
$ws.Range("L3").Value = 21000
$ws.Range("M3").Value = -9136.571599999999
$ws.Range("N3").Value = -21224

$ws.Range("H6").Value = 1165.1111
$ws.Range("I6").Value = 1165.1111
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3495.3333
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = -3382.3333
$ws.Range("M6").ClearContents()

$ws.Range("H113").Value = 947461.75
$ws.Range("I113").Value = 1684033.8
$ws.Range("J113").Value = 440.64285
$ws.Range("K113").Value = 5052101.4
$ws.Range("L113").Value = 1321.92855
$ws.Range("M113").Value = -5049931.4
$ws.Range("N113").Value = -5661.928550000001

$ws.Range("H121").Value = 1255486.8
$ws.Range("I121").Value = 396.66666
$ws.Range("J121").Value = 1389960.6
$ws.Range("K121").Value = 1189.99998
$ws.Range("L121").Value = 4169881.8
$ws.Range("M121").Value = 120.0000199999999
$ws.Range("N121").Value = -4172501.8

$ws.Range("H131").Value = 1732.7609
$ws.Range("J131").Value = 1695.0919
$ws.Range("L131").Value = 5085.2757
$ws.Range("N131").Value = -15165.2757

$ws.Range("H133").Value = 6916.5
$ws.Range("I133").Value = 1999
$ws.Range("J133").Value = 7900
$ws.Range("K133").Value = 5997
$ws.Range("L133").Value = 23700
$ws.Range("M133").Value = -937
$ws.Range("N133").Value = -33820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4943.316
$ws.Range("I70").Value = 5045.0586
$ws.Range("J70").Value = 4860.952
$ws.Range("K70").Value = 5045.0586
$ws.Range("L70").Value = 4860.952
$ws.Range("M70").Value = -4775.0586
$ws.Range("N70").Value = -5400.952

$ws.Range("H73").Value = 4943.316
$ws.Range("I73").Value = 5045.0586
$ws.Range("J73").Value = 4860.952
$ws.Range("K73").Value = 5045.0586
$ws.Range("L73").Value = 4860.952
$ws.Range("M73").Value = -4109.0586
$ws.Range("N73").Value = -6732.952

$ws.Range("H80").Value = 2399.5
$ws.Range("I80").Value = 2399.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2399.5
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = -1401.5
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 2399.5
$ws.Range("I83").Value = 2399.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 11997.5
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = -7005.5
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 24333.334
$ws.Range("J68").Value = 24333.334
$ws.Range("L68").Value = 24333.334
$ws.Range("N68").Value = -25955.334

$ws.Range("H71").Value = 24333.334
$ws.Range("J71").Value = 24333.334
$ws.Range("L71").Value = 73000.00199999999
$ws.Range("N71").Value = -81112.00199999999

$ws.Range("H76").Value = 23000
$ws.Range("J76").Value = 23000
$ws.Range("L76").Value = 23000
$ws.Range("N76").Value = -23630

$ws.Range("H79").Value = 23000
$ws.Range("J79").Value = 23000
$ws.Range("L79").Value = 23000
$ws.Range("N79").Value = -25184

$ws.Range("H113").Value = 743
$ws.Range("I113").Value = 479.23077
$ws.Range("J113").Value = 1124
$ws.Range("K113").Value = 1437.69231
$ws.Range("L113").Value = 3372
$ws.Range("M113").Value = 732.3076900000001
$ws.Range("N113").Value = -7712
